{"js": "// CIV-17609: update GA documents to display \"Case number\" instead of\n// \"Claim number\" (label text only; the merge fields such as\n// <<caseNumber>> are untouched).\nconst body = context.document.body;\n\nconst results = body.search(\"Claim number\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Case number\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# CIV-17609: update GA documents to display \"Case number\" instead of\n# \"Claim number\" (label text only; the merge fields such as\n# <<caseNumber>> are untouched).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\n    \"Claim number\",\n    $true, $true, $false, $false, $false, $true, 1, $false,\n    \"Case number\", 2\n)\n"}
